$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1008.7273
$ws.Range("J2").Value = 1308.1666
$ws.Range("L2").Value = 1308.1666
$ws.Range("N2").Value = -1534.1666

$ws.Range("H64").Value = 6797.864
$ws.Range("I64").Value = 4888.778
$ws.Range("K64").Value = 4888.778
$ws.Range("M64").Value = -4640.778

$ws.Range("H67").Value = 6797.864
$ws.Range("I67").Value = 4888.778
$ws.Range("K67").Value = 4888.778
$ws.Range("M67").Value = -4030.778

$ws.Range("H76").Value = 7312.375
$ws.Range("I76").Value = 7101.6665
$ws.Range("J76").Value = 7944.5
$ws.Range("K76").Value = 7101.6665
$ws.Range("L76").Value = 7944.5
$ws.Range("M76").Value = -6786.6665
$ws.Range("N76").Value = -8574.5

$ws.Range("H79").Value = 7312.375
$ws.Range("I79").Value = 7101.6665
$ws.Range("J79").Value = 7944.5
$ws.Range("K79").Value = 7101.6665
$ws.Range("L79").Value = 7944.5
$ws.Range("M79").Value = -6009.6665
$ws.Range("N79").Value = -10128.5

$ws.Range("H86").Value = 2550.6
$ws.Range("I86").Value = 2001
$ws.Range("J86").Value = 2786.1428
$ws.Range("K86").Value = 2001
$ws.Range("L86").Value = 2786.1428
$ws.Range("M86").Value = -878
$ws.Range("N86").Value = -5032.1428

$ws.Range("H89").Value = 2550.6
$ws.Range("I89").Value = 2001
$ws.Range("J89").Value = 2786.1428
$ws.Range("K89").Value = 10005
$ws.Range("L89").Value = 13930.714
$ws.Range("M89").Value = -4389
$ws.Range("N89").Value = -25162.714

$ws.Range("H132").Value = 45457410
$ws.Range("I132").Value = 58826736
$ws.Range("K132").Value = 176480208
$ws.Range("M132").Value = -176477678

$ws.Range("H133").Value = 85046.664
$ws.Range("J133").Value = 85046.664
$ws.Range("L133").Value = 85046.664
$ws.Range("N133").Value = -95166.664

$ws.Range("H136").Value = 85498
$ws.Range("J136").Value = 85498
$ws.Range("L136").Value = 85498
$ws.Range("N136").Value = -95698

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3789.8542
$ws.Range("I61").Value = 3979.561
$ws.Range("J61").Value = 2678.7144
$ws.Range("K61").Value = 3979.561
$ws.Range("L61").Value = 2678.7144
$ws.Range("M61").Value = -3767.561
$ws.Range("N61").Value = -3102.7144

$ws.Range("H102").Value = 8341841.5
$ws.Range("I102").Value = 9267602
$ws.Range("K102").Value = 9267602
$ws.Range("M102").Value = -9265980

$ws.Range("H136").Value = 3789.8542
$ws.Range("I136").Value = 3979.561
$ws.Range("J136").Value = 2678.7144
$ws.Range("K136").Value = 11938.683
$ws.Range("L136").Value = 8036.1432
$ws.Range("M136").Value = -9388.683000000001
$ws.Range("N136").Value = -13136.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3578841.5
$ws.Range("I94").Value = 4348790
$ws.Range("J94").Value = 37078.4
$ws.Range("K94").Value = 4348790
$ws.Range("L94").Value = 37078.4
$ws.Range("M94").Value = -4348339
$ws.Range("N94").Value = -37980.4

$ws.Range("H99").Value = 8464460
$ws.Range("I99").Value = 9592355
$ws.Range("K99").Value = 9592355
$ws.Range("M99").Value = -9590857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3731.8113
$ws.Range("I31").Value = 5317.074
$ws.Range("K31").Value = 5317.074
$ws.Range("M31").Value = -5022.074

$ws.Range("H34").Value = 3731.8113
$ws.Range("I34").Value = 5317.074
$ws.Range("K34").Value = 5317.074
$ws.Range("M34").Value = -5115.074

$ws.Range("H134").Value = 24519.195
$ws.Range("I134").Value = 31711.941
$ws.Range("J134").Value = 10133.706
$ws.Range("K134").Value = 95135.823
$ws.Range("L134").Value = 30401.118
$ws.Range("M134").Value = -92600.823
$ws.Range("N134").Value = -35471.118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 173.625
$ws.Range("I12").Value = 289
$ws.Range("J12").Value = 157.14285
$ws.Range("K12").Value = 867
$ws.Range("L12").Value = 471.42855
$ws.Range("M12").Value = -694
$ws.Range("N12").Value = -817.4285500000001

$ws.Range("H47").Value = 736
$ws.Range("J47").Value = 1000
$ws.Range("L47").Value = 3000
$ws.Range("N47").Value = -3862

$ws.Range("H64").Value = 4969
$ws.Range("I64").Value = 4969
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 14907
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -14637
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 4969
$ws.Range("I67").Value = 4969
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 14907
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -13971
$ws.Range("N67").ClearContents()

$ws.Range("H128").Value = 198846.88
$ws.Range("I128").Value = 198846.88
$ws.Range("K128").Value = 596540.64
$ws.Range("M128").Value = -591560.64

$ws.Range("H137").Value = 3807
$ws.Range("J137").Value = 5676.6665
$ws.Range("L137").Value = 17029.9995
$ws.Range("N137").Value = -27229.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1291103.4
$ws.Range("I80").Value = 2722677
$ws.Range("J80").Value = 2687
$ws.Range("K80").Value = 2722677
$ws.Range("L80").Value = 2687
$ws.Range("M80").Value = -2721679
$ws.Range("N80").Value = -4683

$ws.Range("H83").Value = 1291103.4
$ws.Range("I83").Value = 2722677
$ws.Range("J83").Value = 2687
$ws.Range("K83").Value = 13613385
$ws.Range("L83").Value = 13435
$ws.Range("M83").Value = -13608393
$ws.Range("N83").Value = -23419

$ws.Range("H132").Value = 10044.714
$ws.Range("I132").Value = 6913.654
$ws.Range("K132").Value = 20740.962
$ws.Range("M132").Value = -18210.962

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7419.0557
$ws.Range("I40").Value = 5903.1333
$ws.Range("J40").Value = 14998.667
$ws.Range("K40").Value = 5903.1333
$ws.Range("L40").Value = 14998.667
$ws.Range("M40").Value = -5767.1333
$ws.Range("N40").Value = -15270.667

$ws.Range("H93").Value = 10103904
$ws.Range("I93").Value = 14495603
$ws.Range("J93").Value = 2996.2
$ws.Range("K93").Value = 14495603
$ws.Range("L93").Value = 2996.2
$ws.Range("M93").Value = -14494355
$ws.Range("N93").Value = -5492.2

$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 10608.952
$ws.Range("I132").Value = 11341.677
$ws.Range("J132").Value = 7494.875
$ws.Range("K132").Value = 34025.031
$ws.Range("L132").Value = 22484.625
$ws.Range("M132").Value = -31495.031
$ws.Range("N132").Value = -27544.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12821873
$ws.Range("I81").Value = 12821873
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 25643746
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -25642685
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 12821873
$ws.Range("I84").Value = 12821873
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 128218730
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -128213426
$ws.Range("N84").ClearContents()

$ws.Range("H96").Value = 2209.3
$ws.Range("I96").Value = 1946.5
$ws.Range("J96").Value = 2275
$ws.Range("K96").Value = 1946.5
$ws.Range("L96").Value = 2275
$ws.Range("M96").Value = -573.5
$ws.Range("N96").Value = -5021

$ws.Range("H107").Value = 37042936
$ws.Range("I107").Value = 55559040
$ws.Range("J107").Value = 10732.333
$ws.Range("K107").Value = 166677120
$ws.Range("L107").Value = 32196.999
$ws.Range("M107").Value = -166675200
$ws.Range("N107").Value = -36036.999

$ws.Range("H132").Value = 24645038
$ws.Range("I132").Value = 27036190
$ws.Range("K132").Value = 81108570
$ws.Range("M132").Value = -81106040

$ws.Range("H135").Value = 100298.2
$ws.Range("I135").Value = 97999
$ws.Range("J135").Value = 101283.57
$ws.Range("K135").Value = 97999
$ws.Range("L135").Value = 101283.57
$ws.Range("M135").Value = -92929
$ws.Range("N135").Value = -111423.57
